$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values -----------------------------------------------------
# Shared strings must land in the workbook in this exact order:
#   0 hours | 1 date | 2 activity | 3 finding-docs... | 4 making-ppt... | 5 total hours
$ws.Range("B1").Value = "hours"
$ws.Range("A1").Value = "date"
$ws.Range("C1").Value = "activity"
$ws.Range("C2").Value = "finding documents and coming up with a proposal and reviewing "
$ws.Range("C3").Value = "Making the powerpoint, continuing to refine proposal, document reviewal"
$ws.Range("E1").Value = "total hours"

$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1

# Dates as raw serials (Nov 4 2025 / Nov 5 2025) so no implicit date-format
# style is fabricated on assignment -- we set the number format ourselves.
$ws.Range("A2").Value = 45965
$ws.Range("A3").Value = 45966

$ws.Range("E2").Formula = "=SUM(B:B)"

# --- Formatting --------------------------------------------------------
$ws.Range("A1:C3").WrapText = $true
$ws.Range("E1:E2").WrapText = $true

$ws.Range("A2:A3").NumberFormat = "mm-dd-yy"

$ws.Columns("A").ColumnWidth = 8.7
$ws.Columns("C").ColumnWidth = 24.92

# --- Selection -----------------------------------------------------
$ws.Range("E3").Select() | Out-Null
